$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value 0.1.0 -> 0.1.1
$ws.Range("B3").Value = "0.1.1"

# Update Date value
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new row for "Jurisdiction" after the "Contact" row (row 10),
# pushing "Description" and everything below down by one row.
$ws.Rows.Item(11).Insert()

# Copy formatting (borders/alignment) from the row below (now row 12,
# the former "Description" row) onto the newly inserted row 11 so the
# new row matches the table's existing style instead of Excel's default
# insert style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
